# Scheduled-runner style refresh of cached market-price / profit values
# across the Sargatanas_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values are plain cached numbers (no formulas in the sheet), so each changed
# cell is written directly with its new value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 223.70589
$ws.Range("J53").Value = 258.2
$ws.Range("L53").Value = 258.2
$ws.Range("N53").Value = -1532.2

# Row 62
$ws.Range("H62").Value = 76958370
$ws.Range("I62").Value = 142858620
$ws.Range("J62").Value = 74734.336
$ws.Range("K62").Value = 142858620
$ws.Range("L62").Value = 74734.336
$ws.Range("M62").Value = -142857996
$ws.Range("N62").Value = -75982.336

# Row 65
$ws.Range("H65").Value = 76958370
$ws.Range("I65").Value = 142858620
$ws.Range("J65").Value = 74734.336
$ws.Range("K65").Value = 714293100
$ws.Range("L65").Value = 373671.68
$ws.Range("M65").Value = -714289980
$ws.Range("N65").Value = -379911.68

# Row 69
$ws.Range("H69").Value = 4993.5
$ws.Range("I69").Value = 4993.5
$ws.Range("K69").Value = 14980.5
$ws.Range("M69").Value = -14106.5

# Row 72
$ws.Range("H72").Value = 4993.5
$ws.Range("I72").Value = 4993.5
$ws.Range("K72").Value = 44941.5
$ws.Range("M72").Value = -40573.5

# Row 76
$ws.Range("H76").Value = 14606.615
$ws.Range("J76").Value = 9216.333000000001
$ws.Range("L76").Value = 9216.333000000001
$ws.Range("N76").Value = -9846.333000000001

# Row 79
$ws.Range("H79").Value = 14606.615
$ws.Range("J79").Value = 9216.333000000001
$ws.Range("L79").Value = 9216.333000000001
$ws.Range("N79").Value = -11400.333

# Row 80
$ws.Range("H80").Value = 25759.55
$ws.Range("I80").Value = 10268.1
$ws.Range("J80").Value = 41251
$ws.Range("K80").Value = 30804.3
$ws.Range("L80").Value = 123753
$ws.Range("M80").Value = -29806.3
$ws.Range("N80").Value = -125749

# Row 83
$ws.Range("H83").Value = 25759.55
$ws.Range("I83").Value = 10268.1
$ws.Range("J83").Value = 41251
$ws.Range("K83").Value = 92412.90000000001
$ws.Range("L83").Value = 371259
$ws.Range("M83").Value = -87420.90000000001
$ws.Range("N83").Value = -381243

# Row 98
$ws.Range("H98").Value = 5434.3945
$ws.Range("I98").Value = 5152.0347
$ws.Range("K98").Value = 5152.0347
$ws.Range("M98").Value = -3654.0347

# Row 103
$ws.Range("H103").Value = 648.5789
$ws.Range("J103").Value = 661.625
$ws.Range("L103").Value = 1984.875
$ws.Range("N103").Value = -3156.875

# Row 122
$ws.Range("H122").Value = 5434.3945
$ws.Range("I122").Value = 5152.0347
$ws.Range("K122").Value = 15456.1041
$ws.Range("M122").Value = -13006.1041

# Row 132
$ws.Range("H132").Value = 2675.7917
$ws.Range("I132").Value = 2675.7917
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8027.375100000001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5497.375100000001
$ws.Range("N132").ClearContents()

# Row 133
$ws.Range("H133").Value = 98983.336
$ws.Range("J133").Value = 98983.336
$ws.Range("L133").Value = 98983.336
$ws.Range("N133").Value = -109103.336

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 4173184
$ws.Range("I97").Value = 484.66666
$ws.Range("K97").Value = 484.66666
$ws.Range("M97").Value = 11.33334000000002

# Row 102
$ws.Range("H102").Value = 1926.1765
$ws.Range("I102").Value = 1926.1765
$ws.Range("K102").Value = 1926.1765
$ws.Range("M102").Value = -304.1765

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 607.11536
$ws.Range("I94").Value = 245.65218
$ws.Range("J94").Value = 3378.3333
$ws.Range("K94").Value = 245.65218
$ws.Range("L94").Value = 3378.3333
$ws.Range("M94").Value = 205.34782
$ws.Range("N94").Value = -4280.3333

# Row 99
$ws.Range("H99").Value = 4787613.5
$ws.Range("I99").Value = 3003.9412
$ws.Range("K99").Value = 3003.9412
$ws.Range("M99").Value = -1505.9412

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 535
$ws.Range("I22").Value = 330
$ws.Range("K22").Value = 330
$ws.Range("M22").Value = 20

# Row 31
$ws.Range("H31").Value = 8709.046
$ws.Range("I31").Value = 3756.35
$ws.Range("J31").Value = 12836.292
$ws.Range("K31").Value = 3756.35
$ws.Range("L31").Value = 12836.292
$ws.Range("M31").Value = -3461.35
$ws.Range("N31").Value = -13426.292

# Row 34
$ws.Range("H34").Value = 8709.046
$ws.Range("I34").Value = 3756.35
$ws.Range("J34").Value = 12836.292
$ws.Range("K34").Value = 3756.35
$ws.Range("L34").Value = 12836.292
$ws.Range("M34").Value = -3554.35
$ws.Range("N34").Value = -13240.292

# Row 58
$ws.Range("H58").Value = 11911041
$ws.Range("I58").Value = 27779814
$ws.Range("J58").Value = 9460
$ws.Range("K58").Value = 27779814
$ws.Range("L58").Value = 9460
$ws.Range("M58").Value = -27779611
$ws.Range("N58").Value = -9866

# Row 62
$ws.Range("H62").Value = 11367408
$ws.Range("I62").Value = 20836540
$ws.Range("J62").Value = 4449.8
$ws.Range("K62").Value = 20836540
$ws.Range("L62").Value = 4449.8
$ws.Range("M62").Value = -20835916
$ws.Range("N62").Value = -5697.8

# Row 65
$ws.Range("H65").Value = 11367408
$ws.Range("I65").Value = 20836540
$ws.Range("J65").Value = 4449.8
$ws.Range("K65").Value = 104182700
$ws.Range("L65").Value = 22249
$ws.Range("M65").Value = -104179580
$ws.Range("N65").Value = -28489

# Row 94
$ws.Range("H94").Value = 986.4545000000001
$ws.Range("I94").Value = 1367.4
$ws.Range("J94").Value = 669
$ws.Range("K94").Value = 1367.4
$ws.Range("L94").Value = 669
$ws.Range("M94").Value = -916.4000000000001
$ws.Range("N94").Value = -1571

# Row 132
$ws.Range("H132").Value = 5252.5
$ws.Range("I132").Value = 3056.2222
$ws.Range("K132").Value = 9168.6666
$ws.Range("M132").Value = -6638.6666

# Row 136
$ws.Range("H136").Value = 11911041
$ws.Range("I136").Value = 27779814
$ws.Range("J136").Value = 9460
$ws.Range("K136").Value = 83339442
$ws.Range("L136").Value = 28380
$ws.Range("M136").Value = -83336892
$ws.Range("N136").Value = -33480

$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 11905319
$ws.Range("I14").Value = 11905319
$ws.Range("K14").Value = 35715957
$ws.Range("M14").Value = -35715784

# Row 92
$ws.Range("H92").Value = 12821290
$ws.Range("I92").Value = 393
$ws.Range("J92").Value = 25642188
$ws.Range("K92").Value = 1179
$ws.Range("L92").Value = 76926564
$ws.Range("M92").Value = 69
$ws.Range("N92").Value = -76929060

# Row 132
$ws.Range("H132").Value = 7412.8237
$ws.Range("I132").Value = 4863.4707
$ws.Range("K132").Value = 43771.2363
$ws.Range("M132").Value = -41241.2363

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 8336.714
$ws.Range("I70").Value = 7306.4165
$ws.Range("K70").Value = 7306.4165
$ws.Range("M70").Value = -7036.4165

# Row 73
$ws.Range("H73").Value = 8336.714
$ws.Range("I73").Value = 7306.4165
$ws.Range("K73").Value = 7306.4165
$ws.Range("M73").Value = -6370.4165

# Row 80
$ws.Range("H80").Value = 51952.3
$ws.Range("I80").Value = 1748.1333
$ws.Range("K80").Value = 1748.1333
$ws.Range("M80").Value = -750.1333

# Row 83
$ws.Range("H83").Value = 51952.3
$ws.Range("I83").Value = 1748.1333
$ws.Range("K83").Value = 8740.666499999999
$ws.Range("M83").Value = -3748.666499999999

# Row 97
$ws.Range("H97").Value = 1624.88
$ws.Range("I97").Value = 1330.6
$ws.Range("J97").Value = 2802
$ws.Range("K97").Value = 1330.6
$ws.Range("L97").Value = 2802
$ws.Range("M97").Value = -834.5999999999999
$ws.Range("N97").Value = -3794

# Row 102
$ws.Range("H102").Value = 712.0769
$ws.Range("I102").Value = 604.75
$ws.Range("K102").Value = 604.75
$ws.Range("M102").Value = 1017.25

# Row 113
$ws.Range("H113").Value = 5802.364
$ws.Range("I113").Value = 3787.8635
$ws.Range("J113").Value = 7816.864
$ws.Range("K113").Value = 3787.8635
$ws.Range("L113").Value = 7816.864
$ws.Range("M113").Value = -1617.8635
$ws.Range("N113").Value = -12156.864

# Row 122
$ws.Range("H122").Value = 2873.2559
$ws.Range("I122").Value = 2311.7742
$ws.Range("J122").Value = 4323.75
$ws.Range("K122").Value = 6935.3226
$ws.Range("L122").Value = 12971.25
$ws.Range("M122").Value = -4485.3226
$ws.Range("N122").Value = -17871.25

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 4023.1924
$ws.Range("I61").Value = 2742.2632
$ws.Range("K61").Value = 2742.2632
$ws.Range("M61").Value = -2540.2632

# Row 113
$ws.Range("H113").Value = 4023.1924
$ws.Range("I113").Value = 2742.2632
$ws.Range("K113").Value = 2742.2632
$ws.Range("M113").Value = -572.2631999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 3351666.8
$ws.Range("I5").Value = 27500
$ws.Range("K5").Value = 27500
$ws.Range("M5").Value = -27388
